$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# U2 and V2 should now both contain the text "No training occurred"
# instead of the previous numeric values (0 and 100 respectively).
$ws.Range("U2").Value = "No training occurred"
$ws.Range("V2").Value = "No training occurred"
